$d = $word.ActiveDocument

# The previous commit added the "SPRINT BACKLOG 2" heading and its table
# skeleton, but left the per-task-group details blank. This change fills in
# the missing USER STORY / OWNER / ESTIMATED EFFORT(HRS) / STATUS cells for
# the two task groups that make up that sprint backlog:
#   S.NO 1 -> "Implementasi Adopt page"   (owner Prana, 10 hrs, Done)
#   S.NO 2 -> "Implementasi Payment page" (owner Prana, 10 hrs, Done)

# Locate the table that immediately follows the "SPRINT BACKLOG 2" heading
# (rather than hard-coding a table index) so the edit still lands correctly
# even if earlier tables in the document were added to/removed.
$heading = $d.Content
$headingFound = $heading.Find.Execute("SPRINT BACKLOG 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$sprintBacklog2 = $null
if ($headingFound) {
    $heading.Collapse(0)
    for ($i = 1; $i -le $d.Tables.Count; $i++) {
        $candidate = $d.Tables.Item($i)
        if ($candidate.Range.Start -ge $heading.Start) {
            $sprintBacklog2 = $candidate
            break
        }
    }
}
if ($null -eq $sprintBacklog2) {
    # Fallback: the Sprint Backlog 2 table is the 4th table in the document.
    $sprintBacklog2 = $d.Tables.Item(4)
}

# --- Row 2 (S.NO = 1) : Implementasi Adopt page ---
$sprintBacklog2.Cell(2, 2).Range.Text = "Implementasi Adopt page"
$sprintBacklog2.Cell(2, 4).Range.Text = "Prana"
$sprintBacklog2.Cell(2, 5).Range.Text = "10"
$sprintBacklog2.Cell(2, 6).Range.Text = "Done"

# --- Row 5 (S.NO = 2) : Implementasi Payment page ---
$sprintBacklog2.Cell(5, 2).Range.Text = "Implementasi Payment page"
$sprintBacklog2.Cell(5, 4).Range.Text = "Prana"
$sprintBacklog2.Cell(5, 5).Range.Text = "10"
$sprintBacklog2.Cell(5, 6).Range.Text = "Done"
